$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 9.427210000000001
$ws.Range("H2").Value = 28.28163
$ws.Range("I2").Value = 0.2188083857550241
$ws.Range("J2").Value = 0.2188083857550241
$ws.Range("M2").Value = 3.087329333333333
$ws.Range("N2").Value = 9.261987999999999
$ws.Range("O2").Value = 0.1539049749041678
$ws.Range("P2").Value = 0.1539049749041678
$ws.Range("Q2").Value = 29.10490196449333
$ws.Range("R2").Value = 261.9441176804399
$ws.Range("S2").Value = 0.03367569911844845
$ws.Range("T2").Value = 0.03367569911844845
$ws.Range("G3").Value = 9.427210000000001
$ws.Range("H3").Value = 28.28163
$ws.Range("I3").Value = 0.2188083857550241
$ws.Range("J3").Value = 0.2188083857550241
$ws.Range("O3").Value = 0.2832552948356705
$ws.Range("P3").Value = 0.2832552948356705
$ws.Range("Q3").Value = 53.56628395052999
$ws.Range("R3").Value = 482.09655555477
$ws.Range("S3").Value = 0.06197863381955649
$ws.Range("T3").Value = 0.06197863381955648
$ws.Range("G4").Value = 9.427210000000001
$ws.Range("H4").Value = 28.28163
$ws.Range("I4").Value = 0.2188083857550241
$ws.Range("J4").Value = 0.2188083857550241
$ws.Range("M4").Value = 4.823431
$ws.Range("N4").Value = 14.470293
$ws.Range("O4").Value = 0.2404505470122564
$ws.Range("P4").Value = 0.2404505470122564
$ws.Range("Q4").Value = 45.47149695751001
$ws.Range("R4").Value = 409.24347261759
$ws.Range("S4").Value = 0.05261259604566437
$ws.Range("T4").Value = 0.05261259604566436
$ws.Range("G5").Value = 9.427210000000001
$ws.Range("H5").Value = 28.28163
$ws.Range("I5").Value = 0.2188083857550241
$ws.Range("J5").Value = 0.2188083857550241
$ws.Range("M5").Value = 6.467117666666667
$ws.Range("N5").Value = 19.401353
$ws.Range("O5").Value = 0.3223891832479054
$ws.Range("P5").Value = 0.3223891832479053
$ws.Range("Q5").Value = 60.96687633837667
$ws.Range("R5").Value = 548.70188704539
$ws.Range("S5").Value = 0.07054145677135484
$ws.Range("T5").Value = 0.07054145677135483
$ws.Range("I6").Value = 0.3808887290954196
$ws.Range("J6").Value = 0.3808887290954196
$ws.Range("M6").Value = 3.087329333333333
$ws.Range("N6").Value = 9.261987999999999
$ws.Range("O6").Value = 0.1539049749041678
$ws.Range("P6").Value = 0.1539049749041678
$ws.Range("Q6").Value = 50.66409626600932
$ws.Range("R6").Value = 455.9768663940839
$ws.Range("S6").Value = 0.05862067029271091
$ws.Range("T6").Value = 0.05862067029271091
$ws.Range("I7").Value = 0.3808887290954196
$ws.Range("J7").Value = 0.3808887290954196
$ws.Range("O7").Value = 0.2832552948356705
$ws.Range("P7").Value = 0.2832552948356705
$ws.Range("S7").Value = 0.1078887492595069
$ws.Range("T7").Value = 0.1078887492595069
$ws.Range("I8").Value = 0.3808887290954196
$ws.Range("J8").Value = 0.3808887290954196
$ws.Range("M8").Value = 4.823431
$ws.Range("N8").Value = 14.470293
$ws.Range("O8").Value = 0.2404505470122564
$ws.Range("P8").Value = 0.2404505470122564
$ws.Range("Q8").Value = 79.154099265661
$ws.Range("R8").Value = 712.386893390949
$ws.Range("S8").Value = 0.09158490326179679
$ws.Range("T8").Value = 0.09158490326179679
$ws.Range("I9").Value = 0.3808887290954196
$ws.Range("J9").Value = 0.3808887290954196
$ws.Range("M9").Value = 6.467117666666667
$ws.Range("N9").Value = 19.401353
$ws.Range("O9").Value = 0.3223891832479054
$ws.Range("P9").Value = 0.3223891832479053
$ws.Range("Q9").Value = 106.1275415259477
$ws.Range("R9").Value = 955.147873733529
$ws.Range("S9").Value = 0.122794406281405
$ws.Range("T9").Value = 0.122794406281405
$ws.Range("G10").Value = 7.213061
$ws.Range("H10").Value = 21.639183
$ws.Range("I10").Value = 0.1674173200514808
$ws.Range("J10").Value = 0.1674173200514808
$ws.Range("M10").Value = 3.087329333333333
$ws.Range("N10").Value = 9.261987999999999
$ws.Range("O10").Value = 0.1539049749041678
$ws.Range("P10").Value = 0.1539049749041678
$ws.Range("Q10").Value = 22.26909480842266
$ws.Range("R10").Value = 200.421853275804
$ws.Range("S10").Value = 0.02576635844104617
$ws.Range("T10").Value = 0.02576635844104617
$ws.Range("G11").Value = 7.213061
$ws.Range("H11").Value = 21.639183
$ws.Range("I11").Value = 0.1674173200514808
$ws.Range("J11").Value = 0.1674173200514808
$ws.Range("O11").Value = 0.2832552948356705
$ws.Range("P11").Value = 0.2832552948356705
$ws.Range("Q11").Value = 40.98528341667299
$ws.Range("R11").Value = 368.867550750057
$ws.Range("S11").Value = 0.04742184235178
$ws.Range("T11").Value = 0.04742184235177999
$ws.Range("G12").Value = 7.213061
$ws.Range("H12").Value = 21.639183
$ws.Range("I12").Value = 0.1674173200514808
$ws.Range("J12").Value = 0.1674173200514808
$ws.Range("M12").Value = 4.823431
$ws.Range("N12").Value = 14.470293
$ws.Range("O12").Value = 0.2404505470122564
$ws.Range("P12").Value = 0.2404505470122564
$ws.Range("Q12").Value = 34.791702032291
$ws.Range("R12").Value = 313.125318290619
$ws.Range("S12").Value = 0.04025558618570456
$ws.Range("T12").Value = 0.04025558618570455
$ws.Range("G13").Value = 7.213061
$ws.Range("H13").Value = 21.639183
$ws.Range("I13").Value = 0.1674173200514808
$ws.Range("J13").Value = 0.1674173200514808
$ws.Range("M13").Value = 6.467117666666667
$ws.Range("N13").Value = 19.401353
$ws.Range("O13").Value = 0.3223891832479054
$ws.Range("P13").Value = 0.3223891832479053
$ws.Range("Q13").Value = 46.64771422384433
$ws.Range("R13").Value = 419.829428014599
$ws.Range("S13").Value = 0.05397353307295006
$ws.Range("T13").Value = 0.05397353307295005
$ws.Range("G14").Value = 10.03371566666667
$ws.Range("H14").Value = 30.101147
$ws.Range("I14").Value = 0.2328855650980756
$ws.Range("J14").Value = 0.2328855650980756
$ws.Range("M14").Value = 3.087329333333333
$ws.Range("N14").Value = 9.261987999999999
$ws.Range("O14").Value = 0.1539049749041678
$ws.Range("P14").Value = 0.1539049749041678
$ws.Range("Q14").Value = 30.97738470002622
$ws.Range("R14").Value = 278.796462300236
$ws.Range("S14").Value = 0.03584224705196225
$ws.Range("T14").Value = 0.03584224705196225
$ws.Range("G15").Value = 10.03371566666667
$ws.Range("H15").Value = 30.101147
$ws.Range("I15").Value = 0.2328855650980756
$ws.Range("J15").Value = 0.2328855650980756
$ws.Range("O15").Value = 0.2832552948356705
$ws.Range("P15").Value = 0.2832552948356705
$ws.Range("Q15").Value = 57.012505553557
$ws.Range("R15").Value = 513.112549982013
$ws.Range("S15").Value = 0.06596606940482715
$ws.Range("T15").Value = 0.06596606940482713
$ws.Range("G16").Value = 10.03371566666667
$ws.Range("H16").Value = 30.101147
$ws.Range("I16").Value = 0.2328855650980756
$ws.Range("J16").Value = 0.2328855650980756
$ws.Range("M16").Value = 4.823431
$ws.Range("N16").Value = 14.470293
$ws.Range("O16").Value = 0.2404505470122564
$ws.Range("P16").Value = 0.2404505470122564
$ws.Range("Q16").Value = 48.39693519178567
$ws.Range("R16").Value = 435.5724167260711
$ws.Range("S16").Value = 0.05599746151909073
$ws.Range("T16").Value = 0.05599746151909073
$ws.Range("G17").Value = 10.03371566666667
$ws.Range("H17").Value = 30.101147
$ws.Range("I17").Value = 0.2328855650980756
$ws.Range("J17").Value = 0.2328855650980756
$ws.Range("M17").Value = 6.467117666666667
$ws.Range("N17").Value = 19.401353
$ws.Range("O17").Value = 0.3223891832479054
$ws.Range("P17").Value = 0.3223891832479053
$ws.Range("Q17").Value = 64.88921985021011
$ws.Range("R17").Value = 584.0029786518911
$ws.Range("S17").Value = 0.07507978712219548
$ws.Range("T17").Value = 0.07507978712219548
